$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A5").Value = "monkey"
$ws.Range("A5").Select()
